$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(43, 'Cecilia Romero', 6, 1),
    @(44, 'Manuel Fernandez', 6, 1),
    @(48, 'María Aguilera', 8, 3),
    @(51, 'Manuel Rodríguez', 6, 1),
    @(52, 'Pedro Torres', 6, 1),
    @(58, 'Julieta Ramos', 11, 6),
    @(60, 'Cecilia Peralta', 9, 4),
    @(66, 'Sofia Martinez', 8, 3),
    @(76, 'Carlos Sánchez', 6, 1),
    @(78, 'Ana Aguilera', 7, 2),
    @(84, 'Juan Santos', 10, 5),
    @(95, 'Julieta Romero', 12, 7),
    @(100, 'Manuel Torres', 10, 5),
    @(102, 'Ana Pérez', 8, 3),
    @(104, 'Luisa Torres', 9, 4),
    @(105, 'Cecilia Pérez', 12, 7),
    @(107, 'Julieta Fernandez', 11, 6),
    @(113, 'Manuel Carril', 8, 3),
    @(118, 'Luisa Fernandez', 11, 6),
    @(119, 'Manuel Chavez', 11, 6),
    @(122, 'Manuel González', 11, 6),
    @(134, 'Diego Santos', 11, 6),
    @(139, 'Hugo Fernandez', 7, 2),
    @(148, 'Sofia Rodríguez', 12, 7),
    @(149, 'Carlos Romero', 7, 2),
    @(152, 'Julieta Peralta', 11, 6),
    @(154, 'Julieta Sánchez', 12, 7),
    @(157, 'Cecilia Aguilera', 11, 6),
    @(158, 'Julieta Pérez', 6, 1),
    @(159, 'Manuel Carril', 9, 4),
    @(160, 'Carlos Martinez', 7, 2),
    @(161, 'María Fernandez', 8, 3),
    @(171, 'Manuel Pérez', 6, 1),
    @(173, 'María Gómez', 8, 3),
    @(182, 'Julieta Santos', 9, 4),
    @(186, 'Juan Chavez', 11, 6),
    @(188, 'María Peralta', 6, 1),
    @(189, 'Carlos Pérez', 7, 2),
    @(190, 'Julieta Chavez', 9, 4),
    @(198, 'Manuel Chavez', 8, 3),
    @(36, 'Luisa Aguilera', 9, 4),
    @(41, 'Victoria Rodríguez', 6, 1),
    @(42, 'Juan Rodríguez', 6, 1),
    @(49, 'Manuel Torres', 6, 1),
    @(61, 'Cecilia Chavez', 6, 1),
    @(68, 'Luisa González', 6, 1),
    @(91, 'Hugo González', 11, 6),
    @(103, 'Fernando Martinez', 11, 6),
    @(106, 'Carlos Santos', 9, 4),
    @(109, 'Manuel Aguilera', 6, 1),
    @(115, 'Diego Santos', 9, 4),
    @(116, 'Manuel Aguilera', 9, 4),
    @(117, 'Fernando Peralta', 7, 2),
    @(120, 'Pedro Vargas', 6, 1),
    @(124, 'Luisa Romero', 11, 6),
    @(125, 'Diego Santos', 9, 4),
    @(128, 'Manuel Peralta', 6, 1),
    @(130, 'Hugo Sánchez', 6, 1),
    @(133, 'Pedro Martinez', 11, 6),
    @(135, 'Diego González', 7, 2),
    @(137, 'Hugo Gómez', 6, 1),
    @(138, 'Manuel Gómez', 6, 1),
    @(140, 'Manuel Vargas', 6, 1),
    @(141, 'Manuel Santos', 6, 1),
    @(143, 'Carlos Torres', 7, 2),
    @(144, 'Diego Gómez', 6, 1),
    @(153, 'Victoria Chavez', 9, 4),
    @(162, 'Sofia Carril', 9, 4),
    @(163, 'Julieta Martinez', 6, 1),
    @(166, 'Sofia Aguilera', 7, 2),
    @(167, 'María Rodríguez', 7, 2),
    @(170, 'Julieta Peralta', 9, 4),
    @(174, 'Sofia Martinez', 9, 4),
    @(176, 'Pedro Rodríguez', 7, 2),
    @(177, 'Julieta Pérez', 6, 1),
    @(184, 'Manuel González', 9, 4),
    @(185, 'Luisa Peralta', 9, 4)
)

$r = 13
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

Write-Output "Updated $($data.Count) rows, last row = $($r - 1)"